# "Add files via upload" - correction of two typos in the acceptance-criteria
# sheet (HU-criterios de aceptacion 3):
#   G13: "vizualisan" -> "vizualizan"
#   G15: "liquidao"   -> "liquidado"
# and leaving the selection on the last-edited cell (H15) as the author did
# before re-uploading the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("G13").Value = "Se deberá mostrar el precio total en la tabla donde se vizualizan los vehiculos registrados "
$ws.Range("G15").Value = "Se deberá dar un estado de inactivo si el vehiculo ha salido o se ha liquidado "

$ws.Range("H15").Select()
